# Append new exchange-rate log entries (2018-03-27 13:33/13:40/13:41) to the
# allData_sheet and to each per-currency sheet (CNY, USD, GBP, EUR, RUB).

$wb = $excel.ActiveWorkbook

# --- allData_sheet: append 15 new rows (3 timestamps x 5 currencies) ---
$allData = $wb.Worksheets.Item("allData_sheet")

$newLogRows = @(
    @("CNY", "0.0592", "2018-03-27 13:33:00"),
    @("USD", "0.0095", "2018-03-27 13:33:00"),
    @("GBP", "0.0066", "2018-03-27 13:33:00"),
    @("EUR", "0.0076", "2018-03-27 13:33:00"),
    @("RUB", "0.5358", "2018-03-27 13:33:00"),
    @("CNY", "0.0592", "2018-03-27 13:40:00"),
    @("USD", "0.0095", "2018-03-27 13:40:00"),
    @("GBP", "0.0066", "2018-03-27 13:40:00"),
    @("EUR", "0.0076", "2018-03-27 13:40:00"),
    @("RUB", "0.5357", "2018-03-27 13:40:00"),
    @("CNY", "0.0592", "2018-03-27 13:41:00"),
    @("USD", "0.0095", "2018-03-27 13:41:00"),
    @("GBP", "0.0066", "2018-03-27 13:41:00"),
    @("EUR", "0.0076", "2018-03-27 13:41:00"),
    @("RUB", "0.5357", "2018-03-27 13:41:00")
)

$startRow = 7
for ($i = 0; $i -lt $newLogRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newLogRows[$i]
    $rng = $allData.Range("A" + $row + ":C" + $row)
    $rng.NumberFormat = "@"
    $allData.Cells.Item($row, 1).Value = $values[0]
    $allData.Cells.Item($row, 2).Value = $values[1]
    $allData.Cells.Item($row, 3).Value = $values[2]
}

# --- per-currency sheets: append the 3 new (rate, date) rows each ---
$currencySheets = @{
    "CNY" = @(
        @("0.0592", "2018-03-27 13:33:00"),
        @("0.0592", "2018-03-27 13:40:00"),
        @("0.0592", "2018-03-27 13:41:00")
    )
    "USD" = @(
        @("0.0095", "2018-03-27 13:33:00"),
        @("0.0095", "2018-03-27 13:40:00"),
        @("0.0095", "2018-03-27 13:41:00")
    )
    "GBP" = @(
        @("0.0066", "2018-03-27 13:33:00"),
        @("0.0066", "2018-03-27 13:40:00"),
        @("0.0066", "2018-03-27 13:41:00")
    )
    "EUR" = @(
        @("0.0076", "2018-03-27 13:33:00"),
        @("0.0076", "2018-03-27 13:40:00"),
        @("0.0076", "2018-03-27 13:41:00")
    )
    "RUB" = @(
        @("0.5358", "2018-03-27 13:33:00"),
        @("0.5357", "2018-03-27 13:40:00"),
        @("0.5357", "2018-03-27 13:41:00")
    )
}

foreach ($name in @("CNY", "USD", "GBP", "EUR", "RUB")) {
    $ws = $wb.Worksheets.Item($name)
    $rows = $currencySheets[$name]
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $row = 3 + $i
        $values = $rows[$i]
        $rng = $ws.Range("A" + $row + ":B" + $row)
        $rng.NumberFormat = "@"
        $ws.Cells.Item($row, 1).Value = $values[0]
        $ws.Cells.Item($row, 2).Value = $values[1]
    }
}
